$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 256375
$ws.Range("I86").Value = 335166.66
$ws.Range("K86").Value = 335166.66
$ws.Range("M86").Value = -334043.66

$ws.Range("H89").Value = 256375
$ws.Range("I89").Value = 335166.66
$ws.Range("K89").Value = 1675833.3
$ws.Range("M89").Value = -1670217.3

$ws.Range("H109").Value = 33999.332
$ws.Range("J109").Value = 33999.332
$ws.Range("L109").Value = 33999.332
$ws.Range("N109").Value = -36773.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 679.2381
$ws.Range("I74").Value = 654.1177
$ws.Range("K74").Value = 654.1177
$ws.Range("M74").Value = 219.8823

$ws.Range("H77").Value = 679.2381
$ws.Range("I77").Value = 654.1177
$ws.Range("K77").Value = 3270.5885
$ws.Range("M77").Value = 1097.4115

$ws.Range("H81").Value = 44998
$ws.Range("J81").Value = 44998
$ws.Range("L81").Value = 44998
$ws.Range("N81").Value = -46994

$ws.Range("H84").Value = 44998
$ws.Range("J84").Value = 44998
$ws.Range("L84").Value = 134994
$ws.Range("N84").Value = -144978

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3386.625
$ws.Range("I134").Value = 3050.4546
$ws.Range("K134").Value = 9151.363799999999
$ws.Range("M134").Value = -6616.363799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1349.4
$ws.Range("I16").Value = 856.2857
$ws.Range("K16").Value = 856.2857
$ws.Range("M16").Value = -569.2857

$ws.Range("H31").Value = 23708.107
$ws.Range("I31").Value = 1294.125
$ws.Range("J31").Value = 36828.49
$ws.Range("K31").Value = 1294.125
$ws.Range("L31").Value = 36828.49
$ws.Range("M31").Value = -999.125
$ws.Range("N31").Value = -37418.49

$ws.Range("H34").Value = 23708.107
$ws.Range("I34").Value = 1294.125
$ws.Range("J34").Value = 36828.49
$ws.Range("K34").Value = 1294.125
$ws.Range("L34").Value = 36828.49
$ws.Range("M34").Value = -1092.125
$ws.Range("N34").Value = -37232.49

$ws.Range("H58").Value = 2272.2856
$ws.Range("I58").Value = 1844.4445
$ws.Range("J58").Value = 3042.4
$ws.Range("K58").Value = 1844.4445
$ws.Range("L58").Value = 3042.4
$ws.Range("M58").Value = -1641.4445
$ws.Range("N58").Value = -3448.4

$ws.Range("H113").Value = 1349.4
$ws.Range("I113").Value = 856.2857
$ws.Range("K113").Value = 856.2857
$ws.Range("M113").Value = 1313.7143

$ws.Range("H132").Value = 4978.364
$ws.Range("I132").Value = 5979.1113
$ws.Range("J132").Value = 4285.5386
$ws.Range("K132").Value = 17937.3339
$ws.Range("L132").Value = 12856.6158
$ws.Range("M132").Value = -15407.3339
$ws.Range("N132").Value = -17916.6158

$ws.Range("H136").Value = 2272.2856
$ws.Range("I136").Value = 1844.4445
$ws.Range("J136").Value = 3042.4
$ws.Range("K136").Value = 5533.333500000001
$ws.Range("L136").Value = 9127.200000000001
$ws.Range("M136").Value = -2983.333500000001
$ws.Range("N136").Value = -14227.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1414.6666
$ws.Range("I5").Value = 1153.5385
$ws.Range("J5").Value = 1657.1428
$ws.Range("K5").Value = 3460.6155
$ws.Range("L5").Value = 4971.428400000001
$ws.Range("M5").Value = -3348.6155
$ws.Range("N5").Value = -5195.428400000001

$ws.Range("H86").Value = 480.8
$ws.Range("I86").Value = 350
$ws.Range("J86").Value = 495.33334
$ws.Range("K86").Value = 1050
$ws.Range("L86").Value = 1486.00002
$ws.Range("M86").Value = 136
$ws.Range("N86").Value = -3858.00002

$ws.Range("H89").Value = 480.8
$ws.Range("I89").Value = 350
$ws.Range("J89").Value = 495.33334
$ws.Range("K89").Value = 3150
$ws.Range("L89").Value = 4458.00006
$ws.Range("M89").Value = 2778
$ws.Range("N89").Value = -16314.00006

$ws.Range("H107").Value = 496944.78
$ws.Range("I107").Value = 525.2
$ws.Range("J107").Value = 1069736.6
$ws.Range("K107").Value = 1575.6
$ws.Range("L107").Value = 3209209.8
$ws.Range("M107").Value = 344.3999999999999
$ws.Range("N107").Value = -3213049.8

$ws.Range("H135").Value = 1414.6666
$ws.Range("I135").Value = 1153.5385
$ws.Range("J135").Value = 1657.1428
$ws.Range("K135").Value = 10381.8465
$ws.Range("L135").Value = 14914.2852
$ws.Range("M135").Value = -7846.846500000001
$ws.Range("N135").Value = -19984.2852

$ws.Range("H137").Value = 76026.8
$ws.Range("I137").Value = 145401.42
$ws.Range("J137").Value = 15324
$ws.Range("K137").Value = 436204.26
$ws.Range("L137").Value = 45972
$ws.Range("M137").Value = -431104.26
$ws.Range("N137").Value = -56172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 631513
$ws.Range("I107").Value = 289
$ws.Range("J107").Value = 1262737
$ws.Range("K107").Value = 289
$ws.Range("L107").Value = 1262737
$ws.Range("M107").Value = 1631
$ws.Range("N107").Value = -1266577

$ws.Range("H113").Value = 1425.8667
$ws.Range("I113").Value = 383.33334
$ws.Range("J113").Value = 1686.5
$ws.Range("K113").Value = 383.33334
$ws.Range("L113").Value = 1686.5
$ws.Range("M113").Value = 1786.66666
$ws.Range("N113").Value = -6026.5

$ws.Range("H126").Value = 3309.3809
$ws.Range("I126").Value = 3020.2
$ws.Range("J126").Value = 4032.3333
$ws.Range("K126").Value = 9060.599999999999
$ws.Range("L126").Value = 12096.9999
$ws.Range("M126").Value = -6590.599999999999
$ws.Range("N126").Value = -17036.9999

$ws.Range("H132").Value = 2507.8965
$ws.Range("I132").Value = 2220.4614
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6661.3842
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4131.3842
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1735.1111
$ws.Range("I61").Value = 1566.4286
$ws.Range("J61").Value = 1842.4546
$ws.Range("K61").Value = 1566.4286
$ws.Range("L61").Value = 1842.4546
$ws.Range("M61").Value = -1364.4286
$ws.Range("N61").Value = -2246.4546

$ws.Range("H113").Value = 1735.1111
$ws.Range("I113").Value = 1566.4286
$ws.Range("J113").Value = 1842.4546
$ws.Range("K113").Value = 1566.4286
$ws.Range("L113").Value = 1842.4546
$ws.Range("M113").Value = 603.5714
$ws.Range("N113").Value = -6182.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 34000
$ws.Range("J112").Value = 34000
$ws.Range("L112").Value = 34000
$ws.Range("N112").Value = -36954
